# Updated symbol list on Sun Jan 22 04:46:13 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns for the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D2").Value = "300.58"
$ws.Range("E2").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E2").Value = "-0.83%"
$ws.Range("D3").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D3").Value = "37.45"
$ws.Range("E3").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E3").Value = "6.81%"
$ws.Range("D4").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D4").Value = "4.971"
$ws.Range("E4").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E4").Value = "-3.63%"
$ws.Range("D5").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D5").Value = "0.07739"
$ws.Range("E5").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E5").Value = "-0.46%"
$ws.Range("D6").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D6").Value = "2.202"
$ws.Range("E6").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E6").Value = "-7.29%"
$ws.Range("D7").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D7").Value = "7.999"
$ws.Range("D8").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D8").Value = "3.992"
$ws.Range("E8").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E8").Value = "1.19%"
$ws.Range("D9").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D9").Value = "0.9140"
$ws.Range("E9").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E9").Value = "-1.94%"
$ws.Range("D10").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D10").Value = "0.09374"
$ws.Range("E10").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E10").Value = "-5.07%"
$ws.Range("E11").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E11").Value = "0.03%"
$ws.Range("D12").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D12").Value = "0.08448"
$ws.Range("D13").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D13").Value = "0.03534"
$ws.Range("E13").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E13").Value = "6.21%"
$ws.Range("D14").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D14").Value = "0.09925"
$ws.Range("E14").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E14").Value = "-0.04%"
$ws.Range("D15").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D15").Value = "0.001467"
$ws.Range("E15").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E15").Value = "-2.25%"
$ws.Range("D16").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D16").Value = "0.005705"
$ws.Range("E16").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E16").Value = "-1.05%"
$ws.Range("D17").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D17").Value = "3.471"
$ws.Range("E17").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E17").Value = "0.26%"
$ws.Range("E18").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E18").Value = "-4.21%"
$ws.Range("E19").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E19").Value = "3.07%"
$ws.Range("E20").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E20").Value = "-1.41%"
$ws.Range("D21").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D21").Value = "4.551"
$ws.Range("E21").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E21").Value = "6.19%"
$ws.Range("D22").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D22").Value = "0.2227"
$ws.Range("E22").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E22").Value = "-3.12%"
$ws.Range("D23").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D23").Value = "0.04641"
$ws.Range("E23").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E23").Value = "1.72%"
$ws.Range("E24").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E24").Value = "1.13%"
$ws.Range("D25").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D25").Value = "0.004444"
$ws.Range("E25").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E25").Value = "1.74%"
$ws.Range("D26").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D26").Value = "0.0001297"
$ws.Range("E26").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E26").Value = "-0.18%"
$ws.Range("D27").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D27").Value = "0.0004742"
$ws.Range("E27").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E27").Value = "39.76%"
$ws.Range("D39").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D39").Value = "0.01751"
$ws.Range("E39").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E39").Value = "-1.94%"
$ws.Range("D40").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D40").Value = "0.04684"
$ws.Range("E40").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E40").Value = "-2.40%"
$ws.Range("D41").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D41").Value = "0.007864"
$ws.Range("E41").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E41").Value = "1.66%"
$ws.Range("D42").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D42").Value = "0.1388"
$ws.Range("E42").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E42").Value = "-1.86%"
$ws.Range("D43").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D43").Value = "0.007658"
$ws.Range("E43").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E43").Value = "7.59%"
$ws.Range("D44").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D44").Value = "0.002284"
$ws.Range("E44").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E44").Value = "8.91%"
$ws.Range("E45").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E45").Value = "6.61%"
$ws.Range("D46").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D46").Value = "0.00006110"
$ws.Range("E46").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E46").Value = "0.01%"
$ws.Range("E47").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E47").Value = "-0.12%"
$ws.Range("D48").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D48").Value = "8.685"
$ws.Range("E48").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E48").Value = "183.52%"
$ws.Range("E49").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E49").Value = "34.99%"
$ws.Range("E50").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E50").Value = "-0.12%"
$ws.Range("E51").NumberFormat = "@"   # keep as text, not a number
$ws.Range("E51").Value = "-0.12%"
